# Applies the "Completed overhaul of Assign_Subbasin_to_POD.R" changes:
#  - Rename the worksheet
#  - Add Basin_ID / Basin_Num / Grouping columns (AO/AP/AQ) with formulas,
#    for the existing rows (2-13) and for 17 new POD rows (14-30)
#  - Populate the new rows' A/B/AK columns with their new data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet ---
$ws.Name = "Missing_MainStem_GIS"

# --- New header cells (row 1) ---
$ws.Range("AO1").Value = "Basin_ID"
$ws.Range("AP1").Value = "Basin_Num"
$ws.Range("AQ1").Value = "Grouping"

# Match the existing header formatting (bold Calibri 12, grey fill, thick
# bottom border) by copying the format from an existing header cell instead
# of creating a brand-new style.
$ws.Range("A1").Copy()
$ws.Range("AO1:AQ1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# --- Populate new data rows (14-30) : A (index), B (application number), AK (basin) ---
$ws.Range("A14").Value = 1
$ws.Range("B14").Value = "A012919A"
$ws.Range("AK14").Value = "R_03_M"

$ws.Range("A15").Value = 2
$ws.Range("B15").Value = "A013393"
$ws.Range("AK15").Value = "R_24_M"

$ws.Range("A16").Value = 3
$ws.Range("B16").Value = "A015736"
$ws.Range("AK16").Value = "R_18_M"

$ws.Range("A17").Value = 4
$ws.Range("B17").Value = "A015737"
$ws.Range("AK17").Value = "R_18_M"

$ws.Range("A18").Value = 5
$ws.Range("B18").Value = "A019351"
$ws.Range("AK18").Value = "R_14_M"

$ws.Range("A19").Value = 6
$ws.Range("B19").Value = "A026298A"
$ws.Range("AK19").Value = "R_17"

$ws.Range("A20").Value = 7
$ws.Range("B20").Value = "A024929"
$ws.Range("AK20").Value = "R_21_M"

$ws.Range("A21").Value = 8
$ws.Range("B21").Value = "A012919B"
$ws.Range("AK21").Value = "R_03_M"

$ws.Range("A22").Value = 9
$ws.Range("B22").Value = "A013832"
$ws.Range("AK22").Value = "R_10_M"

$ws.Range("A23").Value = 10
$ws.Range("B23").Value = "A016961"
$ws.Range("AK23").Value = "R_10"

$ws.Range("A24").Value = 11
$ws.Range("B24").Value = "A015728B"
$ws.Range("AK24").Value = "R_12_M"

$ws.Range("A25").Value = 12
$ws.Range("B25").Value = "A013217"
$ws.Range("AK25").Value = "R_12_M"

$ws.Range("A26").Value = 13
$ws.Range("B26").Value = "A015779"
$ws.Range("AK26").Value = "R_13_M"

$ws.Range("A27").Value = 14
$ws.Range("B27").Value = "A012510"
$ws.Range("AK27").Value = "R_25"

$ws.Range("A28").Value = 15
$ws.Range("B28").Value = "A020491"
$ws.Range("AK28").Value = "R_25"

$ws.Range("A29").Value = 16
$ws.Range("B29").Value = "A022667"
$ws.Range("AK29").Value = "R_25"

$ws.Range("A30").Value = 17
$ws.Range("B30").Value = "A029070"
$ws.Range("AK30").Value = "R_23"

# --- Basin_Num (AP) literal values for every data row (2-30) ---
$ws.Range("AP2").Value = 6
$ws.Range("AP3").Value = 2
$ws.Range("AP4").Value = 16
$ws.Range("AP5").Value = 4
$ws.Range("AP6").Value = 12
$ws.Range("AP7").Value = 12
$ws.Range("AP8").Value = 2
$ws.Range("AP9").Value = 12
$ws.Range("AP10").Value = 20
$ws.Range("AP11").Value = 9
$ws.Range("AP12").Value = 12
$ws.Range("AP13").Value = 16
$ws.Range("AP14").Value = 3
$ws.Range("AP15").Value = 24
$ws.Range("AP16").Value = 18
$ws.Range("AP17").Value = 18
$ws.Range("AP18").Value = 14
$ws.Range("AP19").Value = 17
$ws.Range("AP20").Value = 21
$ws.Range("AP21").Value = 3
$ws.Range("AP22").Value = 10
$ws.Range("AP23").Value = 10
$ws.Range("AP24").Value = 12
$ws.Range("AP25").Value = 12
$ws.Range("AP26").Value = 13
$ws.Range("AP27").Value = 25
$ws.Range("AP28").Value = 25
$ws.Range("AP29").Value = 25
$ws.Range("AP30").Value = 23

# --- Grouping (AQ) : Upper if Basin_Num < 14, otherwise Lower ---
$ws.Range("AQ2:AQ30").Formula = '=IF(AP2 < 14, "Upper", "Lower")'

# --- Basin_ID (AO) : "U_"/"L_" prefix + Basin_Num, based on Grouping ---
$ws.Range("AO2:AO30").Formula = '=IF(AQ2 = "Upper", "U_", "L_") & AP2'

$wb.Save()
